# Reorders the comma-separated "Recorded By" names in column G.
# For every data row on the active sheet, any multi-value entry in column G
# (e.g. "dnasr281@gmail.com, System") has the order of its comma-separated
# items reversed (e.g. "System, dnasr281@gmail.com"). Single-value cells are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $reversed = $parts[($parts.Length - 1)..0]
        $cell.Value2 = [string]::Join(", ", $reversed)
    }
}
